$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New date-label rows (151-172) need column A written as literal text,
# not auto-converted to a date serial. Format the whole new A-range as
# Text first, write the values, then restore the default "Normal" style
# so no stray formatting/style index is left on the cells.
$dateRange = $ws.Range("A151:A172")
$dateRange.NumberFormat = "@"

$ws.Cells.Item(150, 2).Value = 145000
$ws.Cells.Item(150, 3).Value = 0.75
$ws.Cells.Item(150, 4).Value = 0.75
$ws.Cells.Item(150, 5).Value = 0.75
$ws.Cells.Item(150, 6).Value = 5
$ws.Cells.Item(150, 7).Value = 0.75

$ws.Cells.Item(151, 1).Value = "06-08-2021"
$ws.Cells.Item(151, 2).Value = 105000
$ws.Cells.Item(151, 3).Value = 0.75
$ws.Cells.Item(151, 4).Value = 0.75
$ws.Cells.Item(151, 5).Value = 0.75
$ws.Cells.Item(151, 6).Value = 4
$ws.Cells.Item(151, 7).Value = 0.75

$ws.Cells.Item(152, 1).Value = "09-08-2021"
$ws.Cells.Item(152, 2).Value = 300000
$ws.Cells.Item(152, 3).Value = 0.75
$ws.Cells.Item(152, 4).Value = 0.75
$ws.Cells.Item(152, 5).Value = 0.75
$ws.Cells.Item(152, 6).Value = 4
$ws.Cells.Item(152, 7).Value = 0.75

$ws.Cells.Item(153, 1).Value = "10-08-2021"
$ws.Cells.Item(153, 2).Value = 0
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 0.75

$ws.Cells.Item(154, 1).Value = "11-08-2021"
$ws.Cells.Item(154, 2).Value = 75000
$ws.Cells.Item(154, 3).Value = 0.75
$ws.Cells.Item(154, 4).Value = 0.75
$ws.Cells.Item(154, 5).Value = 0.75
$ws.Cells.Item(154, 6).Value = 4
$ws.Cells.Item(154, 7).Value = 0.75

$ws.Cells.Item(155, 1).Value = "12-08-2021"
$ws.Cells.Item(155, 2).Value = 40000
$ws.Cells.Item(155, 3).Value = 0.75
$ws.Cells.Item(155, 4).Value = 0.75
$ws.Cells.Item(155, 5).Value = 0.75
$ws.Cells.Item(155, 6).Value = 3
$ws.Cells.Item(155, 7).Value = 0.75

$ws.Cells.Item(156, 1).Value = "13-08-2021"
$ws.Cells.Item(156, 2).Value = 120000
$ws.Cells.Item(156, 3).Value = 0.75
$ws.Cells.Item(156, 4).Value = 0.75
$ws.Cells.Item(156, 5).Value = 0.75
$ws.Cells.Item(156, 6).Value = 5
$ws.Cells.Item(156, 7).Value = 0.75

$ws.Cells.Item(157, 1).Value = "16-08-2021"
$ws.Cells.Item(157, 2).Value = 140000
$ws.Cells.Item(157, 3).Value = 0.75
$ws.Cells.Item(157, 4).Value = 0.75
$ws.Cells.Item(157, 5).Value = 0.75
$ws.Cells.Item(157, 6).Value = 6
$ws.Cells.Item(157, 7).Value = 0.75

$ws.Cells.Item(158, 1).Value = "17-08-2021"
$ws.Cells.Item(158, 2).Value = 275000
$ws.Cells.Item(158, 3).Value = 0.75
$ws.Cells.Item(158, 4).Value = 0.75
$ws.Cells.Item(158, 5).Value = 0.75
$ws.Cells.Item(158, 6).Value = 8
$ws.Cells.Item(158, 7).Value = 0.75

$ws.Cells.Item(159, 1).Value = "18-08-2021"
$ws.Cells.Item(159, 2).Value = 55000
$ws.Cells.Item(159, 3).Value = 0.75
$ws.Cells.Item(159, 4).Value = 0.75
$ws.Cells.Item(159, 5).Value = 0.75
$ws.Cells.Item(159, 6).Value = 3
$ws.Cells.Item(159, 7).Value = 0.75

$ws.Cells.Item(160, 1).Value = "19-08-2021"
$ws.Cells.Item(160, 2).Value = 60000
$ws.Cells.Item(160, 3).Value = 0.75
$ws.Cells.Item(160, 4).Value = 0.75
$ws.Cells.Item(160, 5).Value = 0.75
$ws.Cells.Item(160, 6).Value = 3
$ws.Cells.Item(160, 7).Value = 0.75

$ws.Cells.Item(161, 1).Value = "20-08-2021"
$ws.Cells.Item(161, 2).Value = 310000
$ws.Cells.Item(161, 3).Value = 0.75
$ws.Cells.Item(161, 4).Value = 0.75
$ws.Cells.Item(161, 5).Value = 0.75
$ws.Cells.Item(161, 6).Value = 7
$ws.Cells.Item(161, 7).Value = 0.75

$ws.Cells.Item(162, 1).Value = "23-08-2021"
$ws.Cells.Item(162, 2).Value = 25000
$ws.Cells.Item(162, 6).Value = 2
$ws.Cells.Item(162, 7).Value = 0.75

$ws.Cells.Item(163, 1).Value = "24-08-2021"
$ws.Cells.Item(163, 2).Value = 60000
$ws.Cells.Item(163, 3).Value = 0.75
$ws.Cells.Item(163, 4).Value = 0.75
$ws.Cells.Item(163, 5).Value = 0.75
$ws.Cells.Item(163, 6).Value = 3
$ws.Cells.Item(163, 7).Value = 0.75

$ws.Cells.Item(164, 1).Value = "25-08-2021"
$ws.Cells.Item(164, 2).Value = 130000
$ws.Cells.Item(164, 3).Value = 0.75
$ws.Cells.Item(164, 4).Value = 0.75
$ws.Cells.Item(164, 5).Value = 0.75
$ws.Cells.Item(164, 6).Value = 4
$ws.Cells.Item(164, 7).Value = 0.75

$ws.Cells.Item(165, 1).Value = "26-08-2021"
$ws.Cells.Item(165, 2).Value = 160000
$ws.Cells.Item(165, 3).Value = 0.75
$ws.Cells.Item(165, 4).Value = 0.75
$ws.Cells.Item(165, 5).Value = 0.75
$ws.Cells.Item(165, 6).Value = 5
$ws.Cells.Item(165, 7).Value = 0.75

$ws.Cells.Item(166, 1).Value = "27-08-2021"
$ws.Cells.Item(166, 2).Value = 130000
$ws.Cells.Item(166, 3).Value = 0.75
$ws.Cells.Item(166, 4).Value = 0.75
$ws.Cells.Item(166, 5).Value = 0.75
$ws.Cells.Item(166, 6).Value = 5
$ws.Cells.Item(166, 7).Value = 0.75

$ws.Cells.Item(167, 1).Value = "30-08-2021"
$ws.Cells.Item(167, 2).Value = 525000
$ws.Cells.Item(167, 3).Value = 0.75
$ws.Cells.Item(167, 4).Value = 0.75
$ws.Cells.Item(167, 5).Value = 0.75
$ws.Cells.Item(167, 6).Value = 8
$ws.Cells.Item(167, 7).Value = 0.75

$ws.Cells.Item(168, 1).Value = "31-08-2021"
$ws.Cells.Item(168, 2).Value = 0
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0.75

$ws.Cells.Item(169, 1).Value = "01-09-2021"
$ws.Cells.Item(169, 2).Value = 0
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 1.5

$ws.Cells.Item(170, 1).Value = "02-09-2021"
$ws.Cells.Item(170, 2).Value = 0
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 1.5

$ws.Cells.Item(171, 1).Value = "03-09-2021"
$ws.Cells.Item(171, 2).Value = 0
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 1.5

$ws.Cells.Item(172, 1).Value = "06-09-2021"
$ws.Cells.Item(172, 7).Value = 1.5

# Restore default styling on the newly-typed date column so it does not
# keep a custom "Text" number format applied to the cells themselves.
$dateRange.Style = "Normal"

